# PNAD 2009 "agressao" table - correction pass.
#
# The sheet had a two-level header ("sexo", "cor ou raça", "grupos de
# idade", "nível de instrução", "classes de rendimento mensal domiciliar
# per capita") plus a trailing source/footnote block, all stored as
# label-only rows (column A populated, no data in B:G). Those rows are
# removed so the data rows collapse together, and the stray
# "unnamed: 1_level_1" column label in B2 is corrected to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so earlier row numbers stay valid while we work.
$rowsToRemove = @(36, 35, 29, 27, 19, 13, 8, 5)
foreach ($r in $rowsToRemove) {
    $ws.Rows($r).Delete()
}

# Fix the mislabeled column header above "pessoa desconhecida" etc.
$ws.Range("B2").Value2 = "total"
